$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the username value in A2 from "shashi.srinivas" to "supreeth.b"
$ws.Range("A2").Value = "supreeth.b"

# Reset the sheet view: scroll back to show column A (clear topLeftCell="L1")
# and move the active selection to C6 (was T5).
$ws.Range("A1").Select()
$ws.Range("C6").Select()
